# Added 2 new items in Meal Category
# Inserts two new rows (Sinugba Porkchop Bangs and Chicken, Pork Afritada)
# right after the existing "Paksiw na Pata" meal row (row 16), pushing the
# remainder of the menu (Appetizer / Pasta / Dessert / Silog / Noodle
# sections) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 17/18, seeding each from row 16 (copy) so the
# inserted rows inherit row 16's formatting (font/fill/height) before we
# overwrite the values below.
$ws.Rows(16).Copy()
$ws.Rows(17).Insert()
$ws.Rows(16).Copy()
$ws.Rows(17).Insert()

# New row 17: Sinugba Porkchop Bangs and Chicken
$ws.Range("A17").Value = "Sinugba Porkchop Bangs and Chicken"
$ws.Range("B17").Value = "Sinugba is a type of barbecue that is cooked on a grill over charcoal and paired with chicken."
$ws.Range("C17").Value = "Meal"
$ws.Range("D17").Value = "₱700.00"
$ws.Range("E17").Value = "Not specified"

# New row 18: Pork Afritada
$ws.Range("A18").Value = "Pork Afritada"
$ws.Range("B18").Value = "Pork Afritada is a Filipino pork stew that is composed of pork slices along with hotdog, potato and carrot."
$ws.Range("C18").Value = "Meal"
$ws.Range("D18").Value = "₱150.00"
$ws.Range("E18").Value = "Not specified"

# Restore the body-row border (thin right + thin bottom) and exact row
# height on the two new rows; touching cells one at a time (rather than a
# multi-cell range) lets the engine reuse the workbook's existing shared
# cell-style instead of minting a near-duplicate one.
foreach ($r in @(17, 18)) {
    foreach ($col in @("A", "B", "C", "D", "E")) {
        $cell = $ws.Range("$col$r")
        $cell.Borders.Item(10).LineStyle = 1
        $cell.Borders.Item(9).LineStyle = 1
    }
    $ws.Rows($r).RowHeight = 27.36
}
